{"js": "// Update the worksheet date header and every division-problem answer cell.\n// The document is a title paragraph (\"YYYY-MM-DD Weekday\") followed by a\n// table whose non-blank cells each hold one \"a\u00f7b=c, d\" answer string.\n// Replacements are applied strictly in document order (title, then every\n// table cell top-to-bottom/left-to-right) since some \"before\" strings repeat\n// (e.g. \"20\u00f73=6, 2\") but map to different \"after\" strings depending on\n// position.\nconst replacements = [\n  [\"2024-12-23 Monday\", \"2024-12-24 Tuesday\"],\n  [\"80\u00f79=8, 8\", \"79\u00f77=11, 2\"],\n  [\"42\u00f78=5, 2\", \"68\u00f73=22, 2\"],\n  [\"89\u00f72=44, 1\", \"97\u00f76=16, 1\"],\n  [\"27\u00f74=6, 3\", \"25\u00f78=3, 1\"],\n  [\"20\u00f74=5, 0\", \"48\u00f77=6, 6\"],\n  [\"78\u00f73=26, 0\", \"53\u00f74=13, 1\"],\n  [\"35\u00f75=7, 0\", \"69\u00f73=23, 0\"],\n  [\"20\u00f77=2, 6\", \"54\u00f77=7, 5\"],\n  [\"50\u00f74=12, 2\", \"65\u00f77=9, 2\"],\n  [\"32\u00f77=4, 4\", \"62\u00f77=8, 6\"],\n  [\"68\u00f72=34, 0\", \"82\u00f76=13, 4\"],\n  [\"48\u00f77=6, 6\", \"84\u00f73=28, 0\"],\n  [\"20\u00f73=6, 2\", \"32\u00f78=4, 0\"],\n  [\"32\u00f78=4, 0\", \"41\u00f72=20, 1\"],\n  [\"69\u00f74=17, 1\", \"48\u00f72=24, 0\"],\n  [\"50\u00f73=16, 2\", \"83\u00f73=27, 2\"],\n  [\"35\u00f74=8, 3\", \"49\u00f74=12, 1\"],\n  [\"43\u00f72=21, 1\", \"92\u00f79=10, 2\"],\n  [\"48\u00f72=24, 0\", \"62\u00f79=6, 8\"],\n  [\"84\u00f76=14, 0\", \"93\u00f79=10, 3\"],\n  [\"73\u00f78=9, 1\", \"89\u00f73=29, 2\"],\n  [\"21\u00f72=10, 1\", \"98\u00f78=12, 2\"],\n  [\"94\u00f77=13, 3\", \"51\u00f72=25, 1\"],\n  [\"65\u00f79=7, 2\", \"21\u00f78=2, 5\"],\n  [\"20\u00f73=6, 2\", \"27\u00f73=9, 0\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet ptr = 0;\nfor (let i = 0; i < paragraphs.items.length && ptr < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  const [expected, next] = replacements[ptr];\n  if (current === expected) {\n    para.insertText(next, \"Replace\");\n    ptr++;\n  }\n}\nawait context.sync();\n\nif (ptr !== replacements.length) {\n  throw new Error(\n    `Only matched ${ptr} of ${replacements.length} expected paragraphs`\n  );\n}\n", "ps1": "# Update the worksheet date header and every division-problem answer cell.\n# The document is a title paragraph (\"YYYY-MM-DD Weekday\") followed by a\n# single table; only 5 of its 20 rows (1, 5, 9, 13, 17) actually hold text,\n# 5 answer cells each. Cells are addressed by absolute (row, col) so the\n# duplicate \"before\" string (\"20\u00f73=6, 2\", which appears in two different\n# cells with two different replacements) is never ambiguous.\n\n$d = $word.ActiveDocument\n\n# `Range.Text` on a paragraph/cell includes the trailing paragraph mark\n# (chr 13) and, for table cells, the cell mark (chr 7) as well; strip those\n# before comparing against the plain answer strings.\nfunction Strip-Marks([string]$s) {\n    return $s.TrimEnd([char]13, [char]7)\n}\n\n$titleBefore = \"2024-12-23 Monday\"\n$titleAfter  = \"2024-12-24 Tuesday\"\n\n$p = $d.Paragraphs.Item(1)\n$current = Strip-Marks $p.Range.Text\nif ($current -ne $titleBefore) {\n    throw \"Title paragraph mismatch: expected '$titleBefore' but found '$current'\"\n}\n$p.Range.Text = $titleAfter\n\n$cellEdits = @(\n  @{ Row = 1;  Col = 1; Before = \"80\u00f79=8, 8\";   After = \"79\u00f77=11, 2\" },\n  @{ Row = 1;  Col = 2; Before = \"42\u00f78=5, 2\";   After = \"68\u00f73=22, 2\" },\n  @{ Row = 1;  Col = 3; Before = \"89\u00f72=44, 1\";  After = \"97\u00f76=16, 1\" },\n  @{ Row = 1;  Col = 4; Before = \"27\u00f74=6, 3\";   After = \"25\u00f78=3, 1\" },\n  @{ Row = 1;  Col = 5; Before = \"20\u00f74=5, 0\";   After = \"48\u00f77=6, 6\" },\n  @{ Row = 5;  Col = 1; Before = \"78\u00f73=26, 0\";  After = \"53\u00f74=13, 1\" },\n  @{ Row = 5;  Col = 2; Before = \"35\u00f75=7, 0\";   After = \"69\u00f73=23, 0\" },\n  @{ Row = 5;  Col = 3; Before = \"20\u00f77=2, 6\";   After = \"54\u00f77=7, 5\" },\n  @{ Row = 5;  Col = 4; Before = \"50\u00f74=12, 2\";  After = \"65\u00f77=9, 2\" },\n  @{ Row = 5;  Col = 5; Before = \"32\u00f77=4, 4\";   After = \"62\u00f77=8, 6\" },\n  @{ Row = 9;  Col = 1; Before = \"68\u00f72=34, 0\";  After = \"82\u00f76=13, 4\" },\n  @{ Row = 9;  Col = 2; Before = \"48\u00f77=6, 6\";   After = \"84\u00f73=28, 0\" },\n  @{ Row = 9;  Col = 3; Before = \"20\u00f73=6, 2\";   After = \"32\u00f78=4, 0\" },\n  @{ Row = 9;  Col = 4; Before = \"32\u00f78=4, 0\";   After = \"41\u00f72=20, 1\" },\n  @{ Row = 9;  Col = 5; Before = \"69\u00f74=17, 1\";  After = \"48\u00f72=24, 0\" },\n  @{ Row = 13; Col = 1; Before = \"50\u00f73=16, 2\";  After = \"83\u00f73=27, 2\" },\n  @{ Row = 13; Col = 2; Before = \"35\u00f74=8, 3\";   After = \"49\u00f74=12, 1\" },\n  @{ Row = 13; Col = 3; Before = \"43\u00f72=21, 1\";  After = \"92\u00f79=10, 2\" },\n  @{ Row = 13; Col = 4; Before = \"48\u00f72=24, 0\";  After = \"62\u00f79=6, 8\" },\n  @{ Row = 13; Col = 5; Before = \"84\u00f76=14, 0\";  After = \"93\u00f79=10, 3\" },\n  @{ Row = 17; Col = 1; Before = \"73\u00f78=9, 1\";   After = \"89\u00f73=29, 2\" },\n  @{ Row = 17; Col = 2; Before = \"21\u00f72=10, 1\";  After = \"98\u00f78=12, 2\" },\n  @{ Row = 17; Col = 3; Before = \"94\u00f77=13, 3\";  After = \"51\u00f72=25, 1\" },\n  @{ Row = 17; Col = 4; Before = \"65\u00f79=7, 2\";   After = \"21\u00f78=2, 5\" },\n  @{ Row = 17; Col = 5; Before = \"20\u00f73=6, 2\";   After = \"27\u00f73=9, 0\" }\n)\n\n$t = $d.Tables.Item(1)\n\nforeach ($edit in $cellEdits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $cellText = Strip-Marks $cell.Range.Text\n    if ($cellText -ne $edit.Before) {\n        throw \"Cell ($($edit.Row),$($edit.Col)) mismatch: expected '$($edit.Before)' but found '$cellText'\"\n    }\n    $cell.Range.Text = $edit.After\n}\n\nWrite-Output \"Updated title + $($cellEdits.Count) table cells\"\n"}
